$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 724, pushing existing rows 724-756 down to 725-757
$ws.Rows("724:724").Insert()

# Populate the newly inserted row 724 with the new record's data
$ws.Range("A724").Value = 3
$ws.Range("B724").Value = "Femacal de La Calera"
$ws.Range("C724").Value = "Coquimbo"
$ws.Range("D724").Value = 45267
$ws.Range("E724").Value = 5
$ws.Range("F724").Value = "Fruta"
$ws.Range("G724").Value = 100108
$ws.Range("H724").Value = "Tropicales y subtropicales"
$ws.Range("I724").Value = 100108002
$ws.Range("J724").Value = "Mango"
$ws.Range("K724").Value = "Sin especificar"
$ws.Range("L724").Value = "Primera"
$ws.Range("M724").Value = 228
$ws.Range("N724").Value = 10000
$ws.Range("O724").Value = 10000
$ws.Range("P724").Value = 10000
$ws.Range("Q724").Value = "$/bandeja 4 kilos"
$ws.Range("R724").Value = "Perú"
$ws.Range("S724").Value = 2500
$ws.Range("T724").Value = 4
